# Applies the Psychotherapie_Antrag.docx content revision:
# rewrites the case-report paragraphs (sections 1, 2, 3, 4, 5, 6) with the
# revised/anonymised wording, and re-signs the closing line on two lines
# ("Mit freundlichen Gruessen," / "Hans Mueller") separated by a manual line break.

$d = $word.ActiveDocument

# Each entry is the verbatim old paragraph text and its verbatim replacement.
$replacements = @(
    @{ Old = "Martina Muster, 38 Jahre alt, verheiratet, lebt mit ihrem Mann und zwei Töchtern (sechs und drei Jahre) in einem Einfamilienhaus am Stadtrand und arbeitet Teilzeit als Grafikdesignerin in einem Werbeunternehmen. Sie hat keine Geschwister."; New = "Die Patientin, anonymisiert unter dem Code Z290361, ist verheiratet und Mutter zweier Kinder: ein 22-jähriger Sohn, der als Architekt tätig ist, und eine 24-jährige Tochter, die sich noch im Studium befindet. Sie arbeitet 40 Stunden pro Woche als Betreuungsassistentin in Teilzeit in einem Pflegeheim. Die Patientin lebt zusammen mit ihrem Ehemann, der bereits berentet ist, und ihren Kindern in einem Einfamilienhaus. Finanzielle Sorgen bestehen nicht." },
    @{ Old = "Sie berichtet seit etwa sechs Monaten von anhaltender Niedergeschlagenheit und Anhedonie – Spaziergänge oder Treffen mit Freundinnen machen ihr kaum noch Freude. Tagsüber fühlt sie sich oft erschöpft und hat Konzentrationsprobleme, was ihre Kreativität im Job erheblich einschränkt. Einschlaf­störungen und frühes Erwachen sind häufig, verbunden mit Grübeln über berufliche Fehler und Konflikte in der Ehe. Affekt abgeflacht, gelegentliche Stimmungseinbrüche und Wertlosigkeitsgefühle. Beim Erstgespräch wirkte sie motorisch leicht verzögert, mit gedrückter Mimik und zurückhaltender Stimme."; New = "Die Patientin leidet seit mehreren Jahren unter einer Vielzahl psychischer Symptome, darunter Schlafstörungen, innere Unruhe, gereizte Stimmung, Versagens- und Schuldgefühle, vermindertes Selbstwertgefühl, emotionale Instabilität, Energieverlust, Konzentrationsstörungen, Verlust der Lebensfreude und Zukunftsängste. Psychisch erscheint sie wach, bewusstseinsklar und adäquat orientiert, jedoch im Kontakt sehr unsicher und schüchtern. Es besteht ein Grübeln im formalen Gedankengang und Befürchtungen, den alltäglichen Anforderungen nicht gewachsen zu sein. Die Gedächtnis- und Aufmerksamkeitsfunktionen sind unauffällig, und es bestehen keine Hinweise auf Wahrnehmungs- oder Ichstörungen. Die Patientin raucht nicht und konsumiert keinen Alkohol. Suizidalität wird glaubhaft verneint." },
    @{ Old = "Somatisch ist sie ansonsten unauffällig; die bekannte Schilddrüsenunterfunktion ist gut eingestellt, und aktuell nimmt sie kein Psychopharmakon. Ein kardiologisches Konsil im letzten Quartal ergab normalen Blutdruck und unauffälliges EKG. 2023 hatte sie eine vierwöchige Physiotherapie wegen einer Schulterneuralgie, ansonsten keinerlei relevante Vorerkrankungen – keine psychosomatischen Bauchschmerzen oder vegetative Dystonie."; New = "Der somatische Befund und weitere Details sind dem beiliegenden ärztlichen Konsiliarbericht zu entnehmen. Es liegen keine aktuellen psychopharmakologischen Medikationen vor. Die Patientin hat bereits an einer Reha-Maßnahme sowie einer sechswöchigen Behandlung in einer psychosomatischen Klinik teilgenommen." },
    @{ Old = "Martina wuchs als Einzelkind in einer Familie auf, in der ihre Mutter phasenweise unter Angststörungen litt, was zu einer eher distanzierten Beziehung führte. Früh entstand hoher Leistungsdruck in der Schule, sie entwickelte Perfektionismustendenzen und adaptive Coping-Strategien: exzessives Arbeiten oder Rückzug. Seit der Geburt ihrer zweiten Tochter fühlt sie sich zunehmend überfordert und isoliert, der Schichtdienst ihres Mannes erschwert die partnerschaftliche Unterstützung."; New = "Die Patientin wuchs in einem von häuslicher Gewalt und Alkoholmissbrauch geprägten Elternhaus auf. Ihr Vater war häufig impulsiv und unberechenbar, ihre Mutter dominant und fordernd. Die Patientin musste früh Verantwortung im Haushalt übernehmen und fühlte sich von ihren Eltern in ihren Bedürfnissen vernachlässigt. Sie hat zwei Geschwister, zu denen kaum Kontakt besteht. Ihre berufliche Laufbahn begann sie als technische Zeichnerin, eine Entscheidung, die von ihren Eltern getroffen wurde. Später absolvierte sie eine Ausbildung zur Betreuungsassistentin, in der sie bis heute tätig ist." },
    @{ Old = "Diagnostisch ergibt sich eine leichte depressive Episode (ICD-10 F32.0) mit typischen Symptomen von Anhedonie, vermindertem Selbstwert und leichtem Schlafdefizit."; New = "Zum Zeitpunkt der Antragstellung liegt die Diagnose F33.1 vor, eine rezidivierende depressive Störung, gegenwärtig mittelgradige Episode." },
    @{ Old = "Wir haben gemeinsam Therapieziele definiert: Verbesserung der Schlafhygiene durch strukturierte Tagebuchführung und achtsamkeitsbasierte Übungen, Steigerung positiver Aktivitäten (tägliche Spaziergänge, regelmäßige Treffen mit Freundinnen), kognitive Umstrukturierung zur Bearbeitung von Perfektionismus und dysfunktionalen Gedanken sowie Förderung von Ressourcen und Stärkung des Selbstwertgefühls durch klientenzentrierte Interventionen. Für das initiale Setting empfehle ich wöchentliche 50-minütige Sitzungen über mindestens zwölf Einheiten, um Stabilisierung und Symptomreduktion zu erreichen. Eine interdisziplinäre Zusammenarbeit mit dem Hausarzt und gegebenenfalls einer Physiotherapeutin zur Entspannung der Muskelverspannungen ist vorgesehen."; New = "Die Therapieziele umfassen den Aufbau eines guten therapeutischen Bündnisses, die Stabilisierung der Stimmung, die Identifikation und Modifikation dysfunktionaler Kognitionen, die Steigerung des Selbstwertgefühls und die Verbesserung der Abgrenzungsfähigkeit. Geplant sind 36 Stunden Langzeittherapie im wöchentlichen Turnus. Die Prognose ist günstig, da die Patientin hoch motiviert ist, regelmäßig an den Sitzungen teilnimmt und bereits erste positive Veränderungen im Umgang mit ihren Belastungen zeigt." },
    @{ Old = "Mit freundlichen Grüßen, Dr. Sabine Meier"; New = "Mit freundlichen Grüßen,^lHans Müller" }
)

foreach ($r in $replacements) {
    $found = $d.Content.Find.Execute(
        $r.Old, $true, $false, $false, $false, $false,
        $true, 1, $false, $r.New, 2)
    if (-not $found) { throw "Could not find expected paragraph text: $($r.Old.Substring(0, [Math]::Min(60, $r.Old.Length)))" }
}

Write-Output "Applied $($replacements.Count) paragraph replacements."
